$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Perf Eval 1")

$ws2.Range("B2").Value = 1342
$ws2.Range("B3").Value = 173
$ws2.Range("B4").Value = 1787
$ws2.Range("B5").Value = 2711
$ws2.Range("B6").Value = 500
$ws2.Range("B7").Value = 1343
$ws2.Range("B8").Value = 1386

$newChartObj = $ws2.Shapes.AddChart2(201, 51)
$chart = $newChartObj.Chart
$chart.SetSourceData($ws2.Range("A1:B8"))
$ser = $chart.SeriesCollection(1)
Write-Host $ser.Values
Write-Host $ser.Formula
